# Add data for 2022-03-15 (carjacking-by-neighborhood-by-month.xlsx)
# - Rename sheet / header label from "through March 06" to "through March 07"
# - Bump a handful of neighborhood/month counts to reflect the newly-added day

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Sheet tab name: "Through 2022-03-06" -> "Through 2022-03-07"
$ws.Name = "Through 2022-03-07"

# Column header (B1, shared string) label update
$ws.Range("B1").Value = "March 2022 (through March 07)"

# North Lawndale (row 4)
$ws.Range("B4").Value = 3
$ws.Range("W4").Value = 1

# Garfield Park (row 5)
$ws.Range("T5").Value = 1

# Rogers Park (row 6)
$ws.Range("T6").Value = 1

# South Shore (row 7)
$ws.Range("E7").Value = 2

# Chicago Lawn (row 10)
$ws.Range("B10").Value = 3

# Woodlawn (row 13)
$ws.Range("E13").Value = 1

# Humboldt Park (row 15)
$ws.Range("W15").Value = 1

# Little Italy, UIC (row 16)
$ws.Range("W16").Value = 2

# West Pullman (row 21)
$ws.Range("T21").Value = 2

# Grand Crossing (row 26)
$ws.Range("H26").Value = 2

# Belmont Cragin (row 33)
$ws.Range("B33").Value = 2

# Roseland (row 36)
$ws.Range("Q36").Value = 1

# Portage Park (row 75)
$ws.Range("K75").Value = 1

# Streeterville (row 84)
$ws.Range("H84").Value = 1
